$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Change column H values from 60 to 59 for rows 2 through 14
$ws.Range("H2:H14").Value = 59
